# ClueLayout.xlsx update:
# "added rooms, room labels, and door indications to boardpanel
#  updated some locations of labels"
#
# The underlying board-layout grid uses short text codes (room letters,
# optionally suffixed with a marker such as * or # to flag a door/label
# position) in columns B:AB. This pass relocates a handful of those
# room/door-label codes to their correct cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 - room A and room T label relocations
$ws.Range("D5").Value = "A*"
$ws.Range("E5").Value = "A"
$ws.Range("M5").Value = "T*"
$ws.Range("O5").Value = "T"

# Row 6 - room C label relocation
$ws.Range("W6").Value = "C*"
$ws.Range("X6").Value = "C"

# Rows 11-12 - room J door indicator relocation
$ws.Range("X11").Value = "J"
$ws.Range("W12").Value = "J*"

# Rows 16-17 - room N label relocation
$ws.Range("B16").Value = "N*"
$ws.Range("E17").Value = "N"

# Rows 16-17 - room U label relocation
$ws.Range("W16").Value = "U"
$ws.Range("V17").Value = "U*"

# Row 24 - room Q door indicator
$ws.Range("X24").Value = "Q#"

# Rows 26-27 - room S label relocation
$ws.Range("E26").Value = "S"
$ws.Range("C27").Value = "S*"

# Row 26 - room Q label relocation
$ws.Range("V26").Value = "Q*"
$ws.Range("X26").Value = "Q"

# View state: selection moved to Z25 with the sheet scrolled so row 2 is
# the first visible row, and the notes column (AD) widened slightly.
$ws.Range("Z25").Select()
$ws.Columns.Item(30).ColumnWidth = 34.666666666666664
